$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '68.372.34'
$ws.Range("E2").Value = '  +0.13%  '

$ws.Range("D3").Value = '2.651.36'
$ws.Range("E3").Value = '  +0.56%  '

$ws.Range("E4").Value = '  -0.03%  '

$ws.Range("D5").Value = '''597.50'
$ws.Range("E5").Value = '  -0.33%  '

$ws.Range("D6").Value = '''158.49'
$ws.Range("E6").Value = '  +2.76%  '

$ws.Range("E7").Value = '  -0.03%  '

$ws.Range("D8").Value = '''0.544'
$ws.Range("E8").Value = '  -0.35%  '

$ws.Range("E9").Value = '  +3.00%  '

$ws.Range("E10").Value = '  -1.28%  '

$ws.Range("E11").Value = '  +0.78%  '

$ws.Range("E12").Value = '  +1.00%  '

$ws.Range("D13").Value = '''28.15'
$ws.Range("E13").Value = '  +0.75%  '

$ws.Range("D14").Value = '''0.0000190'
$ws.Range("E14").Value = '  +1.32%  '

$ws.Range("D15").Value = '3.129.66'
$ws.Range("E15").Value = '  +0.04%  '

$ws.Range("D16").Value = '68.195.36'
$ws.Range("E16").Value = '  -0.08%  '

$ws.Range("D17").Value = '2.672.20'
$ws.Range("E17").Value = '  +1.49%  '

$ws.Range("D18").Value = '''11.42'
$ws.Range("E18").Value = '  +0.22%  '

$ws.Range("D19").Value = '''364.48'
$ws.Range("E19").Value = '  -0.43%  '

$ws.Range("D20").Value = '''7.35'
$ws.Range("E20").Value = '  -0.78%  '

$ws.Range("E21").Value = '  +4.08%  '

$ws.Range("D22").Value = '''4.82'
$ws.Range("E22").Value = '  -0.44%  '

$ws.Range("E23").Value = '  -2.19%  '

$ws.Range("D24").Value = '''75.32'
$ws.Range("E24").Value = '  +2.24%  '

$ws.Range("E25").Value = '  +0.21%  '

$ws.Range("D26").Value = '''9.76'
$ws.Range("E26").Value = '  -2.81%  '

$ws.Range("D27").Value = '2.784.34'
$ws.Range("E27").Value = '  +0.51%  '

$ws.Range("E28").Value = '  +0.25%  '

$ws.Range("E29").Value = '  +0.08%  '

$ws.Range("D30").Value = '''559.15'
$ws.Range("E30").Value = '  -2.52%  '

$ws.Range("E31").Value = '  +0.37%  '

$ws.Range("E32").Value = '  -0.24%  '

$ws.Range("E33").Value = '  +0.36%  '

$ws.Range("E34").Value = '  -0.81%  '

$ws.Range("D36").Value = '''1.57'
$ws.Range("E36").Value = '  +1.84%  '

$ws.Range("E37").Value = '  +3.14%  '

$ws.Range("D38").Value = '''159.66'
$ws.Range("E38").Value = '  -0.40%  '

$ws.Range("E39").Value = '  +0.99%  '

$ws.Range("E40").Value = '  -2.10%  '

$ws.Range("D41").Value = '''5.35'
$ws.Range("E41").Value = '  -0.34%  '

$ws.Range("D42").Value = '0.0₆0333'
$ws.Range("E42").Value = '  +3.70%  '

$ws.Range("D43").Value = '''2.62'
$ws.Range("E43").Value = '  -0.20%  '

$ws.Range("E44").Value = '  +0.05%  '

$ws.Range("D45").Value = '''158.76'
$ws.Range("E45").Value = '  +1.04%  '

$ws.Range("E46").Value = '  +0.58%  '

$ws.Range("D47").Value = '''22.19'
$ws.Range("E47").Value = '  +1.52%  '

$ws.Range("E48").Value = '  -0.89%  '

$ws.Range("D49").Value = '''0.0779'
$ws.Range("E49").Value = '  +0.06%  '

$ws.Range("D50").Value = '''0.616'
$ws.Range("E50").Value = '  +0.31%  '

$ws.Range("D51").Value = '''0.568'
$ws.Range("E51").Value = '  +1.17%  '
